$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C3: email address, with hyperlink (mailto:) and Hyperlink cell style
$ws.Range("C3").Value = "aa000146@ucf.edu"
[void]$ws.Hyperlinks.Add($ws.Range("C3"), "mailto:aa000146@ucf.edu")
$ws.Range("C3").Style = "Hyperlink"

# D3: number(s) of questions answered incorrectly
$ws.Range("D3").Value = "2, 5"

# F3: list of messages corresponding to the wrong questions
$ws.Range("F3").Value = "m2  `nm5  `n"

# Setting a multi-line value auto-expands the row height; this row should
# stay at the default height (unlike row 2, which was already tall), so
# re-fit it back down.
[void]$ws.Rows.Item(3).AutoFit()

# Update the active selection to C4
[void]$ws.Range("C4").Select()
